# Converts "simple" Word fields (<w:fldSimple w:instr="..."/>) into the
# equivalent "complex" field construct made of begin/instrText/separate/end
# field-character runs, the way Word itself rewrites a field after it has
# been edited by hand. Any other content living in the same paragraph
# (bookmarks, etc.) is preserved as-is.

$d = $word.ActiveDocument

function Get-ParagraphContainingPosition($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $p
        }
    }
    return $null
}

function Convert-FldSimpleToComplexField {
    param([string]$InstrText)

    # Find the (still) simple field carrying this instruction text.
    $target = $null
    foreach ($fld in $d.Fields) {
        if ($fld.Code.Text.Trim() -eq $InstrText) {
            $target = $fld
            break
        }
    }
    if ($target -eq $null) {
        return
    }

    $p = Get-ParagraphContainingPosition($target.Code.Start)
    if ($p -eq $null) {
        return
    }

    # Pull the paragraph's real OOXML so we keep whatever else it contains
    # (bookmarks, run properties, etc.) untouched.
    $bodyXml = $p.Range.WordOpenXML
    $paraInner = $null
    if ($bodyXml -match '(?s)<w:body>(.*?)</w:body>') {
        $bodyInner = $matches[1]
        if ($bodyInner -match '(?s)^<w:p\b[^>]*>(.*?)</w:p>') {
            $paraInner = $matches[1]
        }
    }

    $fieldXml = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
                '<w:r><w:instrText>' + $InstrText + '</w:instrText></w:r>' +
                '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
                '<w:r><w:fldChar w:fldCharType="end"/></w:r>'

    if ($paraInner -ne $null) {
        $fieldPattern = '(?s)<w:fldSimple\b[^>]*?(?:/>|>.*?</w:fldSimple>)'
        $newParaInner = [System.Text.RegularExpressions.Regex]::Replace($paraInner, $fieldPattern, $fieldXml)
    }
    else {
        # Fallback: the paragraph only ever contained the field itself.
        $newParaInner = $fieldXml
    }

    $r = $p.Range
    $r.Collapse(1)
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="2048">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>$newParaInner</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    [void]$r.InsertXML($xml)
}

Convert-FldSimpleToComplexField "m:usercontent zone1"
Convert-FldSimpleToComplexField "m:endusercontent"
